# Updated legacy GSC export data.
#
# The oldest day in the "Chart" coverage table (2025-10-25, row 2) has been
# dropped from the export, shifting every later row up by one. In addition,
# the two most-recent remaining days (now rows 2 and 3) don't yet have final
# "Not indexed" / "Indexed" counts from Google Search Console, so those two
# values are blanked out (their "Impressions" figure is still reported).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete first data row (2025-10-25); everything below shifts up.
$ws.Rows("2:2").Delete()

# The new rows 2-3 (2025-10-26 / 2025-10-27) have not-yet-final "Not indexed"
# / "Indexed" columns — clear them to blank text cells. A leading apostrophe
# forces an explicit (empty) text value rather than simply deleting the
# cells, then ClearFormats drops the quote-prefix formatting it implies so
# the cell style matches the rest of the sheet.
$ws.Range("B2:C3").Value = "'"
$ws.Range("B2:C3").ClearFormats()
